# Generate Report for Handback
#
# For each localized-language sheet (zh-cn, de-de) mark the two tracked
# files as handed back: flip the Status text, stamp "Latest Target File"
# / "Latest Handback File" with the file that was handed back (same as
# the source/handoff file for this round) and record the handback
# timestamp in "Latest Handback DateTime".

$wb = $excel.ActiveWorkbook

# BGR-packed OLE color value matching the workbook's custom hyperlink
# font color (ARGB FF6495ED -> R=0x64 G=0x95 B=0xED).
$hyperlinkColor = 15570276

function Set-HandbackRow($ws, $row, $mdFileName, $mdUrl, $xlfFileName, $xlfUrl, $handbackDateTime) {
    $statusCell = $ws.Cells.Item($row, 2)      # column B - Status
    $targetCell = $ws.Cells.Item($row, 5)      # column E - Latest Target File
    $handbackCell = $ws.Cells.Item($row, 6)    # column F - Latest Handback File
    $dateCell = $ws.Cells.Item($row, 7)        # column G - Latest Handback DateTime

    $statusCell.Value = "Handed back: in sync with en-US"

    $targetCell.Value = $mdFileName
    $targetCell.Font.Underline = $True
    $targetCell.Font.Color = $hyperlinkColor
    $ws.Hyperlinks.Add($targetCell, $mdUrl, "", "", $mdFileName) | Out-Null

    $handbackCell.Value = $xlfFileName
    $handbackCell.Font.Underline = $True
    $handbackCell.Font.Color = $hyperlinkColor
    $ws.Hyperlinks.Add($handbackCell, $xlfUrl, "", "", $xlfFileName) | Out-Null

    $dateCell.Value = $handbackDateTime
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

Set-HandbackRow `
    $wsZhCn `
    2 `
    "1d79beb6-6e58-405c-9564-efece4876858.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/4cc3589d7d32636764dd2268ac81bd2e55d86aee/e2e/1d79beb6-6e58-405c-9564-efece4876858.md" `
    "1d79beb6-6e58-405c-9564-efece4876858.0fa014b19d0b88d001c37e5129ecb1524a4763c1.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d57c6af7a2c3f3049f30a0ad620fba0f14921b54/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/1d79beb6-6e58-405c-9564-efece4876858.0fa014b19d0b88d001c37e5129ecb1524a4763c1.zh-cn.xlf" `
    "2016-03-10 00:53:03"

Set-HandbackRow `
    $wsZhCn `
    3 `
    "94c9d3cb-cd47-4f12-8e86-1c1ed060da00.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/4cc3589d7d32636764dd2268ac81bd2e55d86aee/e2e/94c9d3cb-cd47-4f12-8e86-1c1ed060da00.md" `
    "94c9d3cb-cd47-4f12-8e86-1c1ed060da00.0bdb6a19c127cc2d11eecebb6636aec0a772c7c9.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d57c6af7a2c3f3049f30a0ad620fba0f14921b54/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/94c9d3cb-cd47-4f12-8e86-1c1ed060da00.0bdb6a19c127cc2d11eecebb6636aec0a772c7c9.zh-cn.xlf" `
    "2016-03-10 00:53:03"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

Set-HandbackRow `
    $wsDeDe `
    2 `
    "1d79beb6-6e58-405c-9564-efece4876858.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/4cc3589d7d32636764dd2268ac81bd2e55d86aee/e2e/1d79beb6-6e58-405c-9564-efece4876858.md" `
    "1d79beb6-6e58-405c-9564-efece4876858.0fa014b19d0b88d001c37e5129ecb1524a4763c1.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/717eca5f2be05d842599392767889967843afee2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/1d79beb6-6e58-405c-9564-efece4876858.0fa014b19d0b88d001c37e5129ecb1524a4763c1.de-de.xlf" `
    "2016-03-10 00:53:19"

Set-HandbackRow `
    $wsDeDe `
    3 `
    "94c9d3cb-cd47-4f12-8e86-1c1ed060da00.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/4cc3589d7d32636764dd2268ac81bd2e55d86aee/e2e/94c9d3cb-cd47-4f12-8e86-1c1ed060da00.md" `
    "94c9d3cb-cd47-4f12-8e86-1c1ed060da00.0bdb6a19c127cc2d11eecebb6636aec0a772c7c9.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/717eca5f2be05d842599392767889967843afee2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/94c9d3cb-cd47-4f12-8e86-1c1ed060da00.0bdb6a19c127cc2d11eecebb6636aec0a772c7c9.de-de.xlf" `
    "2016-03-10 00:53:19"
